$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (H1:J1): Costo de ventas / Margen $ / % ---
$ws.Range("H1").Value2 = "COSTO DE VENTAS"
$ws.Range("I1").Value2 = "MARGEN $"
$ws.Range("J1").Value2 = "%"

# H1 header: bold, centered, wrapped, yellow fill
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Interior.Color = 65535
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").WrapText = $true

# I1:J1 headers: bold + centered
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108

# Row 1 taller to fit the wrapped header text
$ws.Rows.Item(1).RowHeight = 30

# --- Column I: Margen $ = Base imp. IVA - Costo de ventas ---
$ws.Range("I2").Formula = "=+D2-H2"
$ws.Range("I3:I9").Formula = "=+D3-H3"
$ws.Range("I2:I9").NumberFormat = "#,##0.00"

# --- Column J: % = Margen $ / Costo de ventas (Percent style) ---
$ws.Range("J2:J9").Style = "Percent"
$ws.Range("J2").Formula = "=+I2/H2"
$ws.Range("J3:J9").Formula = "=+I3/H3"

# --- Totals row 10 ---
$ws.Range("D10").Formula = "=SUM(D2:D9)"
$ws.Range("H10").Formula = "=SUM(H2:H9)"

$ws.Range("D10").Font.Bold = $true
$ws.Range("H10").Font.Bold = $true
$ws.Range("D10").NumberFormat = "#,##0.00"
$ws.Range("D10").Borders.Weight = -4138
$ws.Range("H10").Borders.Weight = -4138

# Bottom rows get a slightly taller, thick-bottom look
$ws.Rows.Item(9).RowHeight = 15.75
$ws.Rows.Item(10).RowHeight = 15.75

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# Restore the author's final selection
$ws.Range("E8").Select()

Write-Host "edit applied"
